$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# The sheet is protected, so unprotect it before editing and
# re-protect it afterwards to preserve the original state.
$ws.Unprotect()

# "tela cad. Cidade" service entry is removed from the "Nicolas" row.
$ws.Range("B11").Value2 = ""

# A new service entry "tela cidade" is added to the "Telas já prontas" list.
$ws.Range("E7").Value2 = "tela cidade"

# Restore sheet protection (matches original sheetProtection settings).
$ws.Protect("CC21")

# Move the active selection to B11, as in the updated workbook.
$ws.Range("B11").Select()
